$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the existing last header cell (J1) onto the new
# header cell (K1) before writing its text, so the new column matches the
# bold/bordered/centered header formatting used by the rest of the row.
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$ws.Range("K1").Value = "intervention_type"

# Fill in the new "intervention_type" column values for each data row.
$ws.Range("K2").Value = "DEVICE"
$ws.Range("K3").Value = "PROCEDURE"
$ws.Range("K4").Value = "DEVICE"
$ws.Range("K5").Value = "DEVICE"
$ws.Range("K6").Value = "OTHER"
$ws.Range("K7").Value = "OTHER"
